$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($r = 71; $r -le 127; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $cell.Value = "ok"
    $cell.HorizontalAlignment = -4108
}
